# Updates the cryptos price/volume table (D & E columns) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.697.77"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "1.591.86"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.85"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.27"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0591"
$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").Value = "1.818.37"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "1.566.60"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  -2.17%  "

$ws.Range("D16").Value = "27.688.23"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.21"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.69"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.77"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.84"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.00"
$ws.Range("E26").Value = "  +4.47%  "

$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.105"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0474"
$ws.Range("E31").Value = "  +1.76%  "

$ws.Range("E32").Value = "  -2.38%  "

$ws.Range("D33").Value = "1.382.37"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").Value = "  +0.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  +0.48%  "

$ws.Range("E38").Value = "  +2.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.534"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.984"
$ws.Range("E42").Value = "  +2.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.44"
$ws.Range("E43").Value = "  +0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.18"

$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "1.729.64"
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.85"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("E51").Value = "  -0.07%  "
